# Add a new feedback response (row 6) to the "Form Responses 1" sheet,
# extending the Form_Responses table and the hidden _FilterDatabase
# defined name to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing data row (row 5) down into a brand new
# row 6 so the new row inherits the exact same cell styles/number
# formats that rows 4-5 already use, then overwrite its values.
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(6).Insert(-4121)   # xlShiftDown
$ws.Rows.Item(6).RowHeight = 22.5

$ws.Range("A6").Value = 45985.74253560185
$ws.Range("B6").Value = "mucabap@gmail.com"
$ws.Range("C6").Value = "Before we started working together, my biggest challenge was feeling stuck in my career development. I had goals and ambitions, but I lacked clarity about the best path forward and the strategies needed to progress with confidence. I also struggled to translate my potential into visible results and didn't have a structured plan to improve my performance and positioning.".Replace("didn't", "didn’t")
$ws.Range("D6").Value = "I had the ambition to transition from the design field to the software area, but I often felt that I lacked the technical direction, confidence, and clarity to make that change effectively. This created a sense of stagnation, as I knew I had the potential to contribute more strategically to the tech environment but wasn't fully sure how to bridge the gap between where I was and where I wanted to be. ".Replace("wasn't", "wasn’t")
$ws.Range("E6").Value = "Instead of treating the shift as a leap into the unknown, you helped me break it down into practical steps: identifying the skills I already had that were transferable, mapping the new capabilities I needed to develop, and creating a realistic learning roadmap. This framework gave me clarity, direction, and confidence, and made the transition feel achievable rather than overwhelming"
$ws.Range("F6").Value = "Your mentoring style stood out because it wasn't generic or theoretical — it was personalized, practical, and grounded in real-world experience. Instead of giving broad advice like many others do, you took the time to understand my goals, strengths, and challenges, and then tailored the guidance to my specific situation.".Replace("wasn't", "wasn’t")
$ws.Range("G6").Value = "I noticed a clear improvement in my productivity and decision-making: started operating with intention, milestones, and measurable progress. This not only accelerated my development but also increased my sense of ownership and confidence in my career evolution."
$ws.Range("H6").Value = "Our work helped me accelerate my transition into the software area, prioritize the right skills, and position myself more strategically in the company. Additionally, your technique on how to use AI the right way has been a game changer for me. Learning how to craft effective prompts and leverage AI thoughtfully has not only improved my technical results but also strengthened my communication, decision-making, and overall professional performance."
$ws.Range("I6").Value = "Mansour is not just about receiving advice — it's about gaining clarity, structure, and a personalized strategy for real professional growth. He understands your goals deeply, challenges you with the right questions, and provides practical frameworks you can actually apply. His mentoring combines empathy, market experience, and a results-focused mindset. He can provide the guidance that accelerates your development, sharpens how you think, and helps you take confident steps toward your career goals.".Replace("it's", "it’s")
$ws.Range("J6").Value = $ws.Range("J4").Value
$ws.Range("K6").Value = $ws.Range("K4").Value

# Grow the table to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K6"))

# Keep the hidden AutoFilter-backing defined name in sync with the
# table's new extent.
$name = $wb.Names.Item("_xlnm._FilterDatabase")
$name.RefersTo = "='Form Responses 1'!`$A`$1:`$K`$6"
